$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 113.3518946666667
$ws.Range("H2").Value = 340.055684
$ws.Range("I2").Value = 0.323929285314747
$ws.Range("J2").Value = 0.323929285314747
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 122.328922
$ws.Range("N2").Value = 366.986766
$ws.Range("O2").Value = 0.9783373008518612
$ws.Range("P2").Value = 0.9783373008518613
$ws.Range("Q2").Value = 13866.21508123088
$ws.Range("R2").Value = 124795.935731078
$ws.Range("S2").Value = 0.3169121026617021
$ws.Range("T2").Value = 0.316912102661702

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 113.3518946666667
$ws.Range("H3").Value = 340.055684
$ws.Range("I3").Value = 0.323929285314747
$ws.Range("J3").Value = 0.323929285314747
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.3863573333333334
$ws.Range("N3").Value = 1.159072
$ws.Range("O3").Value = 0.003089929874945324
$ws.Range("P3").Value = 0.003089929874945324
$ws.Range("Q3").Value = 43.79433575169423
$ws.Range("R3").Value = 394.1490217652481
$ws.Range("S3").Value = 0.001000918776063724
$ws.Range("T3").Value = 0.001000918776063724

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 113.3518946666667
$ws.Range("H4").Value = 340.055684
$ws.Range("I4").Value = 0.323929285314747
$ws.Range("J4").Value = 0.323929285314747
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.322294
$ws.Range("N4").Value = 6.966882000000001
$ws.Range("O4").Value = 0.0185727692731934
$ws.Range("P4").Value = 0.0185727692731934
$ws.Range("Q4").Value = 263.2364248730321
$ws.Range("R4").Value = 2369.127823857289
$ws.Range("S4").Value = 0.006016263876981233
$ws.Range("T4").Value = 0.006016263876981232

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 183.9871776666667
$ws.Range("H5").Value = 551.961533
$ws.Range("I5").Value = 0.5257859618835901
$ws.Range("J5").Value = 0.52578596188359
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 122.328922
$ws.Range("N5").Value = 366.986766
$ws.Range("O5").Value = 0.9783373008518612
$ws.Range("P5").Value = 0.9783373008518613
$ws.Range("Q5").Value = 22506.95310578581
$ws.Range("R5").Value = 202562.5779520723
$ws.Range("S5").Value = 0.5143960187749911
$ws.Range("T5").Value = 0.5143960187749911

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 183.9871776666667
$ws.Range("H6").Value = 551.961533
$ws.Range("I6").Value = 0.5257859618835901
$ws.Range("J6").Value = 0.52578596188359
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.3863573333333334
$ws.Range("N6").Value = 1.159072
$ws.Range("O6").Value = 0.003089929874945324
$ws.Range("P6").Value = 0.003089929874945324
$ws.Range("Q6").Value = 71.08479533081956
$ws.Range("R6").Value = 639.7631579773761
$ws.Range("S6").Value = 0.001624641751450968
$ws.Range("T6").Value = 0.001624641751450968

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 183.9871776666667
$ws.Range("H7").Value = 551.961533
$ws.Range("I7").Value = 0.5257859618835901
$ws.Range("J7").Value = 0.52578596188359
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.322294
$ws.Range("N7").Value = 6.966882000000001
$ws.Range("O7").Value = 0.0185727692731934
$ws.Range("P7").Value = 0.0185727692731934
$ws.Range("Q7").Value = 427.272318772234
$ws.Range("R7").Value = 3845.450868950107
$ws.Range("S7").Value = 0.00976530135714798
$ws.Range("T7").Value = 0.009765301357147978

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 52.588828
$ws.Range("H8").Value = 157.766484
$ws.Range("I8").Value = 0.1502847528016631
$ws.Range("J8").Value = 0.1502847528016631
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 122.328922
$ws.Range("N8").Value = 366.986766
$ws.Range("O8").Value = 0.9783373008518612
$ws.Range("P8").Value = 0.9783373008518613
$ws.Range("Q8").Value = 6433.134638483415
$ws.Range("R8").Value = 57898.21174635074
$ws.Range("S8").Value = 0.1470291794151683
$ws.Range("T8").Value = 0.1470291794151682

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 52.588828
$ws.Range("H9").Value = 157.766484
$ws.Range("I9").Value = 0.1502847528016631
$ws.Range("J9").Value = 0.1502847528016631
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.3863573333333334
$ws.Range("N9").Value = 1.159072
$ws.Range("O9").Value = 0.003089929874945324
$ws.Range("P9").Value = 0.003089929874945324
$ws.Range("Q9").Value = 20.31807934920534
$ws.Range("R9").Value = 182.862714142848
$ws.Range("S9").Value = 0.0004643693474306318
$ws.Range("T9").Value = 0.0004643693474306316

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 52.588828
$ws.Range("H10").Value = 157.766484
$ws.Range("I10").Value = 0.1502847528016631
$ws.Range("J10").Value = 0.1502847528016631
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.322294
$ws.Range("N10").Value = 6.966882000000001
$ws.Range("O10").Value = 0.0185727692731934
$ws.Range("P10").Value = 0.0185727692731934
$ws.Range("Q10").Value = 122.126719731432
$ws.Range("R10").Value = 1099.140477582888
$ws.Range("S10").Value = 0.002791204039064195
$ws.Range("T10").Value = 0.002791204039064194
